$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Exam2 (row 10): fill in Angelica's (column D) grade ---
$ws.Range("D10").Value = 95

# --- Weighted score for Exam2 (row 13): new cell for Angelica ---
$ws.Range("D13").Formula = "=(D10+D12)*D11"

# --- Exam3 (row 14): fill in grades for most students ---
# Angelica (D) got a normal (black) entry; everyone else's was entered
# by the instructor afterwards and is shown in red.
$ws.Range("D14").Value = 98

$ws.Range("C14").Value = 95
$ws.Range("E14").Value = 95
$ws.Range("F14").Value = 95
$ws.Range("G14").Value = 95
$ws.Range("H14").Value = 95
$ws.Range("I14").Value = 95
$ws.Range("J14").Value = 95
$ws.Range("C14").Font.Color = 255
$ws.Range("E14").Font.Color = 255
$ws.Range("F14").Font.Color = 255
$ws.Range("G14").Font.Color = 255
$ws.Range("H14").Font.Color = 255
$ws.Range("I14").Font.Color = 255
$ws.Range("J14").Font.Color = 255

# --- Final Project (row 23): fill in grades for most students ---
$ws.Range("D23").Value = 95

$ws.Range("C23").Value = 94
$ws.Range("E23").Value = 90
$ws.Range("F23").Value = 98
$ws.Range("G23").Value = 96
$ws.Range("H23").Value = 86
$ws.Range("I23").Value = 94
$ws.Range("J23").Value = 96

# --- sum row (31): include the Final Project (25) and Exam3 (16) rows ---
$ws.Range("L31").Formula = "=SUM(L9,L19,L13,L25,L16)"

# --- opt sum row (32): same updated total ---
$ws.Range("L32").Formula = "=SUM(L9,L19,L13,L16,L25)"

# --- Current Grade (row 34): Angelica's percentage, newly computed ---
$ws.Range("D34").Formula = '=D31/$L$31'

# --- restore selection ---
$ws.Range("D28").Select() | Out-Null
